$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-56 from serial 45185 to 45204
$ws.Range("C2:C56").Value = 45204
